$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("A2").Value = 53472
$ws.Range("B2").Value = "Stephany Castro"
$ws.Range("C2").Value = "Atendimento ao Cliente"
$ws.Range("D2").Value = "Problemas pessoais"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 45100
$ws.Range("G2").Value = 4308.64

# Row 3
$ws.Range("A3").Value = 1160
$ws.Range("B3").Value = "Dr. Alexandre Carvalho"
$ws.Range("C3").Value = "Recursos Humanos"
$ws.Range("D3").Value = "Doença"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 45100
$ws.Range("G3").Value = 7400.67

# Row 4
$ws.Range("A4").Value = 77297
$ws.Range("B4").Value = "Cauã da Mata"
$ws.Range("D4").Value = "Consulta médica"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 45094
$ws.Range("G4").Value = 6347.54

# Row 5
$ws.Range("A5").Value = 12816
$ws.Range("B5").Value = "Isabel Freitas"
$ws.Range("C5").Value = "Jurídico"
$ws.Range("D5").Value = "Consulta médica"
$ws.Range("E5").Value = 7
$ws.Range("F5").Value = 45100
$ws.Range("G5").Value = 3639.61

# Row 6
$ws.Range("A6").Value = 68949
$ws.Range("B6").Value = "João Lucas Dias"
$ws.Range("C6").Value = "Vendas"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 45080
$ws.Range("G6").Value = 2652.3

# Row 7
$ws.Range("A7").Value = 32137
$ws.Range("B7").Value = "Marcelo Oliveira"
$ws.Range("C7").Value = "Engenharia"
$ws.Range("D7").Value = "Outros"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = 45097
$ws.Range("G7").Value = 6484.91

# Row 8
$ws.Range("A8").Value = 9786
$ws.Range("B8").Value = "Maria Clara Campos"
$ws.Range("E8").Value = 6
$ws.Range("F8").Value = 45104
$ws.Range("G8").Value = 3950.11

# Row 9
$ws.Range("A9").Value = 4424
$ws.Range("B9").Value = "Isaac Ferreira"
$ws.Range("C9").Value = "P&D"
$ws.Range("D9").Value = "Outros"
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 45101
$ws.Range("G9").Value = 9361.799999999999

# Row 10
$ws.Range("A10").Value = 13059
$ws.Range("B10").Value = "Brenda Nunes"
$ws.Range("C10").Value = "Atendimento ao Cliente"
$ws.Range("D10").Value = "Problemas pessoais"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = 45090
$ws.Range("G10").Value = 6266.46

# Row 11
$ws.Range("A11").Value = 83115
$ws.Range("B11").Value = "Vitória Campos"
$ws.Range("C11").Value = "Financeiro"
$ws.Range("D11").Value = "Viagem de negócios"
$ws.Range("E11").Value = 4
$ws.Range("F11").Value = 45083
$ws.Range("G11").Value = 5115.28
